$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.333.16'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.377.56'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.72'
$ws.Range('E5').Value = '  +2.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.61'
$ws.Range('E6').Value = '  +1.59%  '
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.378.13'
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('E10').Value = '  +3.33%  '
$ws.Range('E11').Value = '  +2.22%  '
$ws.Range('E12').Value = '  +2.19%  '
$ws.Range('E13').Value = '  +3.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.51'
$ws.Range('E14').Value = '  +3.05%  '
$ws.Range('E15').Value = '  +6.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.806.67'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.299.38'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.377.07'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.93'
$ws.Range('E19').Value = '  +3.44%  '
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '320.32'
$ws.Range('E21').Value = '  +1.68%  '
$ws.Range('E22').Value = '  +1.90%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.27'
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('E25').Value = '  -8.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.86'
$ws.Range('E26').Value = '  +4.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.493.96'
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('E29').Value = '  +2.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '519.10'
$ws.Range('E30').Value = '  +3.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0901'
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('E33').Value = '  +2.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.84'
$ws.Range('E34').Value = '  +3.29%  '
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +5.02%  '
$ws.Range('E38').Value = '  +2.79%  '
$ws.Range('E39').Value = '  +6.38%  '
$ws.Range('E40').Value = '  +1.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.51'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '146.32'
$ws.Range('E42').Value = '  +5.58%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.34'
$ws.Range('E44').Value = '  +3.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.80'
$ws.Range('E45').Value = '  +6.73%  '
$ws.Range('E46').Value = '  +2.25%  '
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('E48').Value = '  +2.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.74'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('E50').Value = '  +2.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0908'
$ws.Range('E51').Value = '  +1.44%  '
